$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.380.51"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "'1.627.09"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'0.9994"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'304.39"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").Value = "'0.3786"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'51.91"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "'0.3625"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.228"
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.08099"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'22.73"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "'6.555"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "'0.00001247"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "'7.223"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "'1.624.89"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "'93.57"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'0.06901"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'6.417"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").Value = "'23.368.12"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "'12.71"
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("D25").Value = "'3.256"
$ws.Range("E25").Value = "  +3.93%  "
$ws.Range("D26").Value = "'2.445"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "'149.83"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").Value = "'134.12"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "'2.308"
$ws.Range("E31").Value = "  -4.87%  "
$ws.Range("D32").Value = "'1.806.54"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'6.784"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").Value = "'11.05"
$ws.Range("E34").Value = "  +5.73%  "
$ws.Range("D35").Value = "'0.9519"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "'0.02781"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").Value = "'0.2520"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "'0.08831"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "'6.089"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").Value = "'0.07133"
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").Value = "'1.362"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("D42").Value = "'0.7061"
$ws.Range("D43").Value = "'16.12"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("D45").Value = "'0.6453"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "'0.9987"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "'3.991"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'0.07996"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").Value = "'125.78"
$ws.Range("E51").Value = "  -4.59%  "
